$d = $word.ActiveDocument

# 1) Remove the paragraph "Assume every old person in one facility are served simultaneously."
#    entirely (including its paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Assume every old person in one facility are served simultaneously*") {
        $p.Range.Delete()
        break
    }
}

# 2) Insert an additional empty paragraph right after the paragraph ending in
#    "...client 2 has only 1 old person requiring our service), etc." and before
#    the existing block of empty paragraphs. The new paragraph should be a plain
#    empty paragraph (no list formatting / numbering), matching the other blank
#    paragraphs already present there.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*client 2 requires service 2 only*") {
        $newPara = $p.Range.InsertParagraphAfter()
        $p.Next().Range.ParagraphFormat.Reset()
        $p.Next().Range.ListFormat.RemoveNumbers()
        break
    }
}
